# "Links a Clase 9 en tabla" - add the missing Slides link for Clase 9,
# add the .R taller link for Clase 8, turn "C3 (20%)" into a markdown
# link to homework/c_3.html, and leave the final selection on E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clase 10 row: "C3 (20%)" entrega becomes a link to the homework page.
$ws.Range("F12").Value = "[C3 (20%)](homework/c_3.html)"

# Clase 8 row: material link gains the ".R" taller script link.
$ws.Range("G10").Value = "[Slides](slides/class_8/class_8#1)" + [char]160 + "[.qmd](slides/class_8/class_8.qmd) [.R](slides/class_8/class_8_taller.R)"

# Clase 9 row: it was missing its Slides/material link entirely - add it.
$ws.Range("G11").Value = "[Slides](slides/class_9/class_9#1)" + [char]160 + "[.qmd](slides/class_9/class_9.qmd)"

# Leave the selection where the author ended up after the edit.
$ws.Range("E10").Select()
